# Applies the two substantive changes described by the commit:
#
# 1. The table on slide 5 (the "B1 - TYPES OF FINANCIAL DOCUMENTS" table)
#    is switched from the deck's default table style to the built-in
#    style {D42D49B8-CBB2-42F9-A1CD-ECAF2787AA60}.
#
# 2. The presentation's theme palette is reset from the custom
#    "Integral" / "Red Violet" colors to the stock Office theme colors
#    (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink), which is what the
#    canonical deck ends up using for its (shared) theme palette.

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------
$s5 = $p.Slides.Item(5)
$tableShape = $s5.Shapes.Item(2)
$tbl = $tableShape.Table
$tbl.ApplyStyle("{D42D49B8-CBB2-42F9-A1CD-ECAF2787AA60}")

# --- 2. Theme colors --------------------------------------------------
# RGB() values below are COLORREF-encoded (0x00BBGGRR) equivalents of the
# stock "Office" theme palette:
#   dk1=000000 lt1=FFFFFF dk2=44546A lt2=E7E6E6
#   accent1=5B9BD5 accent2=ED7D31 accent3=A5A5A5 accent4=FFC000
#   accent5=4472C4 accent6=70AD47 hlink=0563C1 folHlink=954F72
$cs = $p.SlideMaster.ColorScheme
$cs.Colors(1).RGB  = 0         # dk1      000000
$cs.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$cs.Colors(3).RGB  = 6968388   # dk2      44546A
$cs.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$cs.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$cs.Colors(6).RGB  = 3243501   # accent2  ED7D31
$cs.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$cs.Colors(8).RGB  = 49407     # accent4  FFC000
$cs.Colors(9).RGB  = 12874308  # accent5  4472C4
$cs.Colors(10).RGB = 4697456   # accent6  70AD47
$cs.Colors(11).RGB = 12673797  # hlink    0563C1
$cs.Colors(12).RGB = 7491477   # folHlink 954F72
